# Update the public EPEX spot / gas / CO2 price workbook.
#
# 1) "Prix Spot" sheet: a new "17-dec" day column is inserted right before
#    the "01-oct." column (i.e. before column ET), shifting every column
#    from ET..FX one place to the right (EU..FY). The new column's header
#    cell gets the same header style as its neighbours, and its 24 data
#    cells (rows 2-25) are filled with the placeholder "-" used for
#    not-yet-published days.
# 2) "Gaz" sheet: one new row is appended with the next day's gas price.
# 3) "CO2" sheet: one new row is appended with the next day's CO2 price.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Prix Spot: insert a new day column before ET (01-oct.)
# ---------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

$wsPrix.Range("ET:ET").Insert()

$wsPrix.Range("ET1").Value = "17-dec"

for ($r = 2; $r -le 25; $r++) {
    $wsPrix.Cells.Item($r, 150).Value = "-"
}

# ---------------------------------------------------------------------
# 2) Gaz: append 2025-12-15 price row
# ---------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

# Use a scratch cell pre-formatted as Text to type the date in as a
# literal string (matching the existing column A cells) rather than
# letting it be auto-recognised as a date serial, then paste the value
# (and only the value) into the new row so no stray number format is
# left behind on the target cell.
$wsGaz.Range("D1").NumberFormat = "@"
$wsGaz.Range("D1").Value = "2025-12-15"
$wsGaz.Range("D1").Copy()
$wsGaz.Range("A180").PasteSpecial(-4163)
$wsGaz.Range("D1").Clear()

$wsGaz.Range("B180").Value = 25.68

# ---------------------------------------------------------------------
# 3) CO2: append 2025-12-15 price row
# ---------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")

$wsCo2.Range("D1").NumberFormat = "@"
$wsCo2.Range("D1").Value = "2025-12-15"
$wsCo2.Range("D1").Copy()
$wsCo2.Range("A180").PasteSpecial(-4163)
$wsCo2.Range("D1").Clear()

$wsCo2.Range("B180").Value = 84.59999999999999
